$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16 (existing rows 16-51 shift down to 17-52)
$ws.Rows("16:16").Insert()

# Populate the new row 16 with the DOB entry
$ws.Range("A16").Value = "DOB"
$ws.Range("D16").Value = "/wlq-res-doc:WildlifeLicenseQueryResults/wlq-res-ext:WildlifeLicenseReport/nc:Person/nc:PersonBirthDate/nc:Date"

# Match the style used by similar "s=8" data-row entries (e.g. row 23/35/45)
$ws.Range("A16:D16").Style = $ws.Range("A23:D23").Style
$ws.Rows("16:16").RowHeight = $ws.Rows("23:23").RowHeight

# Update the frozen-pane top-left cell and active selection to match the new view
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("B17").Select()
